$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 9 de Mayo de 2020 a las 23:04"

# Row 4 - Estados Unidos (values refreshed)
$ws.Range("B4").Value = 1341803
$ws.Range("C4").Value = 20018
$ws.Range("E4").Value = 1029548
$ws.Range("G4").Value = 1280
$ws.Range("H4").Value = 79895

# Row 11 - Brasil (values refreshed)
$ws.Range("B11").Value = 148931
$ws.Range("C11").Value = 3039
$ws.Range("E11").Value = 79528
$ws.Range("G11").Value = 114
$ws.Range("H11").Value = 10106

# Rows 16/17 - India & Peru swap rank (Peru moves to row 16, India to row 17) with refreshed values
$ws.Range("A16").Value = "Peru"
$ws.Range("B16").Value = 65015
$ws.Range("C16").Value = 3168
$ws.Range("D16").Value = 20246
$ws.Range("E16").Value = 42955
$ws.Range("F16").Value = 748
$ws.Range("G16").Value = 100
$ws.Range("H16").Value = 1814

$ws.Range("A17").Value = "India"
$ws.Range("B17").Value = 62808
$ws.Range("C17").Value = 3113
$ws.Range("D17").Value = 19301
$ws.Range("E17").Value = 41406
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 116
$ws.Range("H17").Value = 2101

# Row 37 - Rumania (values refreshed)
$ws.Range("E37").Value = 7280
$ws.Range("G37").Value = 16
$ws.Range("H37").Value = 939

# Rows 62/63 - Afganistan & Ghana swap rank with refreshed values
$ws.Range("A62").Value = "Ghana"
$ws.Range("B62").Value = 4263
$ws.Range("C62").Value = 251
$ws.Range("D62").Value = 378
$ws.Range("E62").Value = 3863
$ws.Range("F62").Value = 5
$ws.Range("G62").Value = 4
$ws.Range("H62").Value = 22

$ws.Range("A63").Value = "Afganistan"
$ws.Range("B63").Value = 4033
$ws.Range("C63").Value = 255
$ws.Range("D63").Value = 502
$ws.Range("E63").Value = 3416
$ws.Range("F63").Value = 7
$ws.Range("G63").Value = 6
$ws.Range("H63").Value = 115

# Row 148 - Suazilandia (values refreshed)
$ws.Range("B148").Value = 163
$ws.Range("C148").Value = 4
$ws.Range("D148").Value = 14
$ws.Range("E148").Value = 147

# Rows 159/160 - Trinidad yTobago & Uganda swap rank with refreshed values
$ws.Range("A159").Value = "Uganda"
$ws.Range("B159").Value = 116
$ws.Range("C159").Value = 15
$ws.Range("D159").Value = 55
$ws.Range("E159").Value = 61
$ws.Range("H159").Value = 0

$ws.Range("A160").Value = "Trinidad yTobago"
$ws.Range("B160").Value = 116
$ws.Range("C160").Value = 0
$ws.Range("D160").Value = 104
$ws.Range("E160").Value = 4
$ws.Range("H160").Value = 8

# Rows 192/193 - Belice & Nueva Caledonia swap rank with refreshed values
$ws.Range("A192").Value = "Nueva Caledonia"
$ws.Range("D192").Value = 18
$ws.Range("H192").Value = 0

$ws.Range("A193").Value = "Belice"
$ws.Range("D193").Value = 16
$ws.Range("H193").Value = 2
